$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 40 (pushes old row 40.."49 down to 41.."50,
# and auto-extends the merged ranges that spanned across the insertion
# point: A34:A40 -> A34:A41, A41:A49 -> A42:A50, B41:B49 -> B42:B50,
# B39:B40 -> B39:B41).
$ws.Rows("40:40").Insert()

# The new row sits "in the middle" of both the A34:A41 ("toolbox API")
# and B39:B41 ("user interface/UI") merged groups, so it should carry
# the same border/font formatting as any other interior row of those
# groups. Pull that formatting from existing interior rows rather than
# guessing raw style indices.
$ws.Range("A6").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B10").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C36:D36").Copy()
$ws.Range("C40:D40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New function entry.
$ws.Range("C40").Value = "figureViewer"
$ws.Range("D40").Value = "图片拼接滚动查看"

# Match the saved selection state from the edit.
$ws.Range("D43").Select()
